$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell updates for existing rows 245-271 (per diff) ---
$ws.Range("D245").Value = 44617
$ws.Range("D246").Value = 44617
$ws.Range("D247").Value = 44264
$ws.Range("J247").Value = 600
$ws.Range("L247").Value = 700
$ws.Range("M247").Value = 650
$ws.Range("P247").Value = 130
$ws.Range("D248").Value = 44264
$ws.Range("I248").Value = "Segunda"
$ws.Range("J248").Value = 300
$ws.Range("K248").Value = 500
$ws.Range("L248").Value = 500
$ws.Range("M248").Value = 500
$ws.Range("P248").Value = 100
$ws.Range("D249").Value = 44504
$ws.Range("J249").Value = 1700
$ws.Range("L249").Value = 650
$ws.Range("M249").Value = 624
$ws.Range("P249").Value = 125
$ws.Range("D250").Value = 44572
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 280
$ws.Range("K250").Value = 600
$ws.Range("L250").Value = 650
$ws.Range("M250").Value = 627
$ws.Range("P250").Value = 125
$ws.Range("D251").Value = 44370
$ws.Range("D252").Value = 44370
$ws.Range("D253").Value = 44385
$ws.Range("J253").Value = 600
$ws.Range("D254").Value = 44385
$ws.Range("J254").Value = 300
$ws.Range("D255").Value = 44236
$ws.Range("J255").Value = 800
$ws.Range("D256").Value = 44236
$ws.Range("J256").Value = 400
$ws.Range("D257").Value = 44229
$ws.Range("J257").Value = 600
$ws.Range("D258").Value = 44229
$ws.Range("J258").Value = 300
$ws.Range("D259").Value = 44299
$ws.Range("D260").Value = 44299
$ws.Range("D261").Value = 44610
$ws.Range("J261").Value = 800
$ws.Range("D262").Value = 44610
$ws.Range("J262").Value = 400
$ws.Range("D263").Value = 44399
$ws.Range("D264").Value = 44399
$ws.Range("D265").Value = 44167
$ws.Range("D266").Value = 44167
$ws.Range("D267").Value = 44390
$ws.Range("J267").Value = 600
$ws.Range("D268").Value = 44390
$ws.Range("J268").Value = 300
$ws.Range("D269").Value = 44285
$ws.Range("J269").Value = 800
$ws.Range("L269").Value = 700
$ws.Range("M269").Value = 650
$ws.Range("O269").Value = "Región Metropolitana"
$ws.Range("P269").Value = 130
$ws.Range("D270").Value = 44285
$ws.Range("I270").Value = "Segunda"
$ws.Range("J270").Value = 400
$ws.Range("K270").Value = 500
$ws.Range("L270").Value = 500
$ws.Range("M270").Value = 500
$ws.Range("P270").Value = 100
$ws.Range("D271").Value = 44498
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 1500
$ws.Range("K271").Value = 600
$ws.Range("L271").Value = 650
$ws.Range("M271").Value = 623
$ws.Range("O271").Value = "Región del Maule"
$ws.Range("P271").Value = 125

# --- New rows 272-273 appended ---
$ws.Range("A272").Value = 11
$ws.Range("B272").Value = "Vega Monumental Concepción"
$ws.Range("C272").Value = "Bíobío"
$ws.Range("D272").Value = 44595
$ws.Range("E272").Value = 8
$ws.Range("F272").Value = 100114014
$ws.Range("G272").Value = "Betarraga"
$ws.Range("H272").Value = "Sin especificar"
$ws.Range("I272").Value = "Primera"
$ws.Range("J272").Value = 800
$ws.Range("K272").Value = 600
$ws.Range("L272").Value = 700
$ws.Range("M272").Value = 650
$ws.Range("N272").Value = "`$/paquete 5 unidades"
$ws.Range("O272").Value = "Región Metropolitana"
$ws.Range("P272").Value = 130
$ws.Range("Q272").Value = 5
$ws.Range("R272").Value = "Hortaliza"

$ws.Range("A273").Value = 11
$ws.Range("B273").Value = "Vega Monumental Concepción"
$ws.Range("C273").Value = "Bíobío"
$ws.Range("D273").Value = 44595
$ws.Range("E273").Value = 8
$ws.Range("F273").Value = 100114014
$ws.Range("G273").Value = "Betarraga"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "Segunda"
$ws.Range("J273").Value = 400
$ws.Range("K273").Value = 500
$ws.Range("L273").Value = 500
$ws.Range("M273").Value = 500
$ws.Range("N273").Value = "`$/paquete 5 unidades"
$ws.Range("O273").Value = "Región Metropolitana"
$ws.Range("P273").Value = 100
$ws.Range("Q273").Value = 5
$ws.Range("R273").Value = "Hortaliza"

# Match date number format for the new rows' D column with existing date cells
$ws.Range("D272").NumberFormat = $ws.Range("D270").NumberFormat
$ws.Range("D273").NumberFormat = $ws.Range("D270").NumberFormat

